$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells to Text format before writing,
# so values like "12.40", "1.000" or "0.00000000124" are preserved as
# literal text instead of being auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.048.78'
$ws.Range("D3").Value = '1.829.18'
$ws.Range("D4").Value = '0.9991'
$ws.Range("D5").Value = '241.64'
$ws.Range("D6").Value = '0.6307'
$ws.Range("D8").Value = '44.62'
$ws.Range("D9").Value = '0.2931'
$ws.Range("D10").Value = '0.07334'
$ws.Range("D11").Value = '22.88'
$ws.Range("D12").Value = '0.07679'
$ws.Range("D13").Value = '1.829.20'
$ws.Range("D14").Value = '4.985'
$ws.Range("D15").Value = '0.6633'
$ws.Range("D16").Value = '82.09'
$ws.Range("D17").Value = '6.062'
$ws.Range("D18").Value = '0.000008653'
$ws.Range("D19").Value = '28.879.19'
$ws.Range("D20").Value = '2.081.87'
$ws.Range("D21").Value = '12.40'
$ws.Range("D22").Value = '224.02'
$ws.Range("D24").Value = '7.138'
$ws.Range("D26").Value = '158.03'
$ws.Range("D27").Value = '8.450'
$ws.Range("D28").Value = '0.1370'
$ws.Range("D29").Value = '17.87'
$ws.Range("D30").Value = '1.503'
$ws.Range("D31").Value = '4.097'
$ws.Range("D32").Value = '1.202'
$ws.Range("D34").Value = '0.05301'
$ws.Range("D35").Value = '0.7405'
$ws.Range("D36").Value = '1.828'
$ws.Range("D37").Value = '1.151'
$ws.Range("D39").Value = '1.295.14'
$ws.Range("D40").Value = '2.743'
$ws.Range("D41").Value = '0.01781'
$ws.Range("D42").Value = '6.333'
$ws.Range("D43").Value = '0.8941'
$ws.Range("D45").Value = '102.56'
$ws.Range("D46").Value = '0.00000000124'
$ws.Range("D47").Value = '1.979.01'
$ws.Range("D48").Value = '0.5138'
$ws.Range("D49").Value = '64.21'
$ws.Range("D50").Value = '1.730'
$ws.Range("D51").Value = '0.05829'

# Remove the temporary Text number format again so no extra cell
# style is left behind on the Price column.
$priceRange.ClearFormats()

# Remaining (non price-risk) cell updates: coin name, link and the
# Volume(1h) percentage column.
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("E6").Value = '  -5.93%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +5.92%  '
$ws.Range("E9").Value = '  -0.23%  '
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("E16").Value = '  -4.44%  '
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("E18").Value = '  +4.13%  '
$ws.Range("E19").Value = '  -1.25%  '
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("E21").Value = '  -1.00%  '
$ws.Range("E22").Value = '  -1.69%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -2.89%  '
$ws.Range("E28").Value = '  -2.45%  '
$ws.Range("E29").Value = '  -0.89%  '
$ws.Range("E30").Value = '  -0.58%  '
$ws.Range("E31").Value = '  -1.49%  '
$ws.Range("E32").Value = '  +0.53%  '
$ws.Range("E33").Value = '  -1.24%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  -1.50%  '
$ws.Range("E36").Value = '  -2.68%  '
$ws.Range("E37").Value = '  +1.33%  '
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("E39").Value = '  -1.90%  '
$ws.Range("E40").Value = '  +0.60%  '
$ws.Range("E41").Value = '  -1.34%  '
$ws.Range("E42").Value = '  +6.10%  '
$ws.Range("E43").Value = '  -2.97%  '
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("E45").Value = '  +0.68%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("E46").Value = '  +2.84%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("E47").Value = '  +0.15%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("E48").Value = '  -0.49%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E49").Value = '  +0.55%  '
$ws.Range("E50").Value = '  -2.54%  '
$ws.Range("E51").Value = '  -1.91%  '
